$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.72904439220488
$ws.Range("C2").Value = 8.513261556760259
$ws.Range("D2").Value = 14.37557280116933
$ws.Range("E2").Value = 15.34236555640216
$ws.Range("G2").Value = 3.671290136683783
$ws.Range("J2").Value = 9.03361154662759
$ws.Range("M2").Value = 18.98649499236267
$ws.Range("N2").Value = 18.70640435242177
$ws.Range("O2").Value = 27.96912167992028
$ws.Range("B3").Value = 15.3141241827853
$ws.Range("C3").Value = 8.138979562866673
$ws.Range("D3").Value = 14.36212761784862
$ws.Range("E3").Value = 15.35916240154982
$ws.Range("G3").Value = 3.67396206409995
$ws.Range("J3").Value = 9.053770850526581
$ws.Range("M3").Value = 18.86686280104447
$ws.Range("N3").Value = 18.7749345714912
$ws.Range("O3").Value = 28.00143219801121
$ws.Range("B4").Value = 15.05682538163829
$ws.Range("C4").Value = 7.901284961013491
$ws.Range("D4").Value = 14.35683326751841
$ws.Range("E4").Value = 15.3722365643375
$ws.Range("G4").Value = 3.675689500497684
$ws.Range("J4").Value = 9.066952788009667
$ws.Range("M4").Value = 18.79666573946472
$ws.Range("N4").Value = 18.81894102833706
$ws.Range("O4").Value = 28.02837001936225
$ws.Range("B5").Value = 14.95150124344966
$ws.Range("C5").Value = 7.802578284188582
$ws.Range("D5").Value = 14.35542214591324
$ws.Range("E5").Value = 15.37825816345913
$ws.Range("G5").Value = 3.676415361090493
$ws.Range("J5").Value = 9.072527117835385
$ws.Range("M5").Value = 18.76890265853485
$ws.Range("N5").Value = 18.83736059884187
$ws.Range("O5").Value = 28.0411277921349
$ws.Range("B6").Value = 14.93398867156896
$ws.Range("C6").Value = 7.786081139520936
$ws.Range("D6").Value = 14.35523295178985
$ws.Range("E6").Value = 15.37929993056426
$ws.Range("G6").Value = 3.676537215552519
$ws.Range("J6").Value = 9.073464978398519
$ws.Range("M6").Value = 18.76434418784648
$ws.Range("N6").Value = 18.84044859031654
$ws.Range("O6").Value = 28.04335358914891
$ws.Range("B7").Value = 15.0554066339766
$ws.Range("C7").Value = 7.899961038886347
$ws.Range("D7").Value = 14.35681121261874
$ws.Range("E7").Value = 15.37231496541878
$ws.Range("G7").Value = 3.675699200871844
$ws.Range("J7").Value = 9.06702714458978
$ws.Range("M7").Value = 18.796287874173
$ws.Range("N7").Value = 18.81918746883509
$ws.Range("O7").Value = 28.02853487327318
$ws.Range("B8").Value = 15.58660267360562
$ws.Range("C8").Value = 8.385922749170822
$ws.Range("D8").Value = 14.37032355253423
$ws.Range("E8").Value = 15.34758386864839
$ws.Range("G8").Value = 3.672193433653303
$ws.Range("J8").Value = 9.040395835822684
$ws.Range("M8").Value = 18.94458427325292
$ws.Range("N8").Value = 18.72963426071376
$ws.Range("O8").Value = 27.9787863262155
$ws.Range("B9").Value = 16.6013442663214
$ws.Range("C9").Value = 9.271159946160482
$ws.Range("D9").Value = 14.42021752137766
$ws.Range("E9").Value = 15.32100735247079
$ws.Range("G9").Value = 3.66600447120437
$ws.Range("J9").Value = 8.994533648146499
$ws.Range("M9").Value = 19.26014411225715
$ws.Range("N9").Value = 18.56924883727973
$ws.Range("O9").Value = 27.93773701784268
$ws.Range("B10").Value = 17.32201079452507
$ws.Range("C10").Value = 9.874446273497338
$ws.Range("D10").Value = 14.47097724834274
$ws.Range("E10").Value = 15.31485770502379
$ws.Range("G10").Value = 3.661870822639445
$ws.Range("J10").Value = 8.964692163360906
$ws.Range("M10").Value = 19.50553173558565
$ws.Range("N10").Value = 18.46059125114192
$ws.Range("O10").Value = 27.9422259369082
$ws.Range("B11").Value = 17.64287908492229
$ws.Range("C11").Value = 10.13772546510164
$ws.Range("D11").Value = 14.49708692397432
$ws.Range("E11").Value = 15.3149633977717
$ws.Range("G11").Value = 3.660079075324513
$ws.Range("J11").Value = 8.951948038300023
$ws.Range("M11").Value = 19.6197583691937
$ws.Range("N11").Value = 18.41313033026164
$ws.Range("O11").Value = 27.95181677483936
$ws.Range("B12").Value = 17.7632640908803
$ws.Range("C12").Value = 10.2357548484651
$ws.Range("D12").Value = 14.50740345719558
$ws.Range("E12").Value = 15.3154203992975
$ws.Range("G12").Value = 3.659413261272557
$ws.Range("J12").Value = 8.947241265194265
$ws.Range("M12").Value = 19.66335630724332
$ws.Range("N12").Value = 18.39543942444007
$ws.Range("O12").Value = 27.95653457674108
$ws.Range("B13").Value = 17.73738866136134
$ws.Range("C13").Value = 10.21471759746915
$ws.Range("D13").Value = 14.50516259170724
$ws.Range("E13").Value = 15.3153034408358
$ws.Range("G13").Value = 3.659556093308194
$ws.Range("J13").Value = 8.948249659013802
$ws.Range("M13").Value = 19.65395196260372
$ws.Range("N13").Value = 18.39923698043761
$ws.Range("O13").Value = 27.95547021802327
$ws.Range("B14").Value = 17.6528063403752
$ws.Range("C14").Value = 10.14582417649313
$ws.Range("D14").Value = 14.49792709513224
$ws.Range("E14").Value = 15.31499264231322
$ws.Range("G14").Value = 3.660024044633022
$ws.Range("J14").Value = 8.951558423187068
$ws.Range("M14").Value = 19.62333850646265
$ws.Range("N14").Value = 18.41166925522075
$ws.Range("O14").Value = 27.9521831477105
$ws.Range("B15").Value = 17.60084790250104
$ws.Range("C15").Value = 10.10340592006354
$ws.Range("D15").Value = 14.49355091161546
$ws.Range("E15").Value = 15.31485655366478
$ws.Range("G15").Value = 3.660312327835038
$ws.Range("J15").Value = 8.95360064385995
$ws.Range("M15").Value = 19.60463059674651
$ws.Range("N15").Value = 18.41932100047979
$ws.Range("O15").Value = 27.95031114378152
$ws.Range("B16").Value = 17.30089069542388
$ws.Range("C16").Value = 9.857010098540858
$ws.Range("D16").Value = 14.46933127243959
$ws.Range("E16").Value = 15.31490917976999
$ws.Range("G16").Value = 3.661989695911323
$ws.Range("J16").Value = 8.965541713367257
$ws.Range("M16").Value = 19.49811641900394
$ws.Range("N16").Value = 18.463732407905
$ws.Range("O16").Value = 27.94175111402667
$ws.Range("B17").Value = 17.11500432448227
$ws.Range("C17").Value = 9.702947256448281
$ws.Range("D17").Value = 14.45524328852451
$ws.Range("E17").Value = 15.3156847556133
$ws.Range("G17").Value = 3.663041368284699
$ws.Range("J17").Value = 8.973079750778549
$ws.Range("M17").Value = 19.43341718437607
$ws.Range("N17").Value = 18.49148035575798
$ws.Range("O17").Value = 27.93843385755411
$ws.Range("B18").Value = 17.00743812676889
$ws.Range("C18").Value = 9.613287190700449
$ws.Range("D18").Value = 14.44742471392733
$ws.Range("E18").Value = 15.31640415229436
$ws.Range("G18").Value = 3.66365461255997
$ws.Range("J18").Value = 8.977493659908426
$ws.Range("M18").Value = 19.39645037329953
$ws.Range("N18").Value = 18.50762556244497
$ws.Range("O18").Value = 27.93723641998344
$ws.Range("B19").Value = 16.97091034671378
$ws.Range("C19").Value = 9.582752101000407
$ws.Range("D19").Value = 14.44482647157184
$ws.Range("E19").Value = 15.31669467992909
$ws.Range("G19").Value = 3.663863682786771
$ws.Range("J19").Value = 8.979001579947345
$ws.Range("M19").Value = 19.38397731087648
$ws.Range("N19").Value = 18.51312392627619
$ws.Range("O19").Value = 27.9369529988616
$ws.Range("B20").Value = 17.13486033904223
$ws.Range("C20").Value = 9.719456386709165
$ws.Range("D20").Value = 14.4567135705445
$ws.Range("E20").Value = 15.3155739110327
$ws.Range("G20").Value = 3.662928552256853
$ws.Range("J20").Value = 8.972269219813496
$ws.Range("M20").Value = 19.44027924083166
$ws.Range("N20").Value = 18.48850737001482
$ws.Range("O20").Value = 27.93871343854563
$ws.Range("B21").Value = 17.67768153526892
$ws.Range("C21").Value = 10.16610559375981
$ws.Range("D21").Value = 14.50004072320755
$ws.Range("E21").Value = 15.31507261973743
$ws.Range("G21").Value = 3.659886252352663
$ws.Range("J21").Value = 8.950583327871611
$ws.Range("M21").Value = 19.6323213660142
$ws.Range("N21").Value = 18.40800996621294
$ws.Range("O21").Value = 27.95311916831117
$ws.Range("B22").Value = 18.02586544270345
$ws.Range("C22").Value = 10.44826718271442
$ws.Range("D22").Value = 14.53085765061456
$ws.Range("E22").Value = 15.3171750818165
$ws.Range("G22").Value = 3.657971820920402
$ws.Range("J22").Value = 8.93710469264723
$ws.Range("M22").Value = 19.7598153896276
$ws.Range("N22").Value = 18.3570404682238
$ws.Range("O22").Value = 27.96886350181885
$ws.Range("B23").Value = 17.84067197147699
$ws.Range("C23").Value = 10.29858292496659
$ws.Range("D23").Value = 14.51418300659386
$ws.Range("E23").Value = 15.31583082502667
$ws.Range("G23").Value = 3.658986850943076
$ws.Range("J23").Value = 8.944235070045202
$ws.Range("M23").Value = 19.69159827666346
$ws.Range("N23").Value = 18.38409425345192
$ws.Range("O23").Value = 27.95988140961406
$ws.Range("B24").Value = 17.12588559867336
$ws.Range("C24").Value = 9.711995990638639
$ws.Range("D24").Value = 14.45604798104806
$ws.Range("E24").Value = 15.31562317187815
$ws.Range("G24").Value = 3.662979529514657
$ws.Range("J24").Value = 8.972635411040965
$ws.Range("M24").Value = 19.43717618742821
$ws.Range("N24").Value = 18.48985085695407
$ws.Range("O24").Value = 27.93858482926774
$ws.Range("B25").Value = 16.33062418132669
$ws.Range("C25").Value = 9.039587768839962
$ws.Range("D25").Value = 14.40422968781613
$ws.Range("E25").Value = 15.32584814079054
$ws.Range("G25").Value = 3.667605815881241
$ws.Range("J25").Value = 9.006262103061218
$ws.Range("M25").Value = 19.17228237256042
$ws.Range("N25").Value = 18.61101806192204
$ws.Range("O25").Value = 27.94276849254996
